# Implemented automatic persistence forecasting for wind
#
# 1. Rename "Input File Settings" sheet to "Input Data Settings".
# 2. On that sheet, rename the "Filename" column header to "Data Source".
# 3. Replace the wind 15-min-ahead forecast filename row with a "persistence"
#    data source entry (automatic persistence forecasting instead of a file).
# 4. Remove the now-unused "Use Persistence Forecast From Actual" column (G).
# 5. Refresh column A's width to fit the new header/content and leave the
#    selection where the author's edit session ended up.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Input File Settings")

$ws2.Name = "Input Data Settings"

# Set A5 (wind forecast row) before A1 (header) so new shared-string entries
# are appended in the same order as the reference edit ("persistence" then
# "Data Source").
$ws2.Range("A5").Value = "persistence"
$ws2.Range("A1").Value = "Data Source"

# The "Use Persistence Forecast From Actual" boolean column is no longer
# needed now that persistence is modeled as its own data source row.
$ws2.Columns.Item(7).Delete()

$ws2.Columns.Item(1).AutoFit()

$ws2.Range("A10").Select()
